$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '23.604.64'
Set-TextValue 'E2' '  +0.67%  '
Set-TextValue 'D3' '1.646.54'
Set-TextValue 'E3' '  +1.02%  '
Set-TextValue 'D4' '0.9981'
Set-TextValue 'E4' '  -0.30%  '
Set-TextValue 'D5' '0.9975'
Set-TextValue 'E5' '  -0.36%  '
Set-TextValue 'D6' '304.34'
Set-TextValue 'E6' '  +0.05%  '
Set-TextValue 'D7' '0.3794'
Set-TextValue 'E7' '  +0.61%  '
Set-TextValue 'D8' '52.02'
Set-TextValue 'E8' '  +0.90%  '
Set-TextValue 'D9' '0.3607'
Set-TextValue 'E9' '  -1.30%  '
Set-TextValue 'D10' '1.244'
Set-TextValue 'E10' '  +1.27%  '
Set-TextValue 'D11' '0.08196'
Set-TextValue 'E11' '  -0.64%  '
Set-TextValue 'D12' '0.9981'
Set-TextValue 'E12' '  -0.37%  '
Set-TextValue 'D13' '22.45'
Set-TextValue 'E13' '  +0.01%  '
Set-TextValue 'D14' '6.530'
Set-TextValue 'E14' '  -0.52%  '
Set-TextValue 'D15' '7.367'
Set-TextValue 'E15' '  +0.75%  '
Set-TextValue 'D16' '0.00001228'
Set-TextValue 'E16' '  -2.00%  '
Set-TextValue 'D17' '1.647.31'
Set-TextValue 'E17' '  +1.13%  '
Set-TextValue 'D18' '97.07'
Set-TextValue 'E18' '  +3.19%  '
Set-TextValue 'D19' '0.06984'
Set-TextValue 'E19' '  +0.05%  '
Set-TextValue 'D20' '6.731'
Set-TextValue 'E20' '  +3.89%  '
Set-TextValue 'D21' '17.54'
Set-TextValue 'E21' '  -1.13%  '
Set-TextValue 'D22' '0.9971'
Set-TextValue 'E22' '  -0.37%  '
Set-TextValue 'D23' '12.57'
Set-TextValue 'E23' '  -1.17%  '
Set-TextValue 'D24' '23.612.95'
Set-TextValue 'E24' '  +0.75%  '
Set-TextValue 'D25' '2.522'
Set-TextValue 'E25' '  +1.95%  '
Set-TextValue 'D26' '3.117'
Set-TextValue 'E26' '  -2.50%  '
Set-TextValue 'D27' '21.27'
Set-TextValue 'E27' '  -0.59%  '
Set-TextValue 'D28' '152.25'
Set-TextValue 'E28' '  +1.61%  '
Set-TextValue 'D29' '5.183'
Set-TextValue 'E29' '  -2.53%  '
Set-TextValue 'D30' '134.76'
Set-TextValue 'E30' '  +0.54%  '
Set-TextValue 'D31' '1.830.31'
Set-TextValue 'E31' '  +1.05%  '
Set-TextValue 'D32' '6.763'
Set-TextValue 'E32' '  -0.57%  '
Set-TextValue 'D33' '1.093'
Set-TextValue 'E33' '  +7.28%  '
Set-TextValue 'D34' '2.047'
Set-TextValue 'E34' '  -9.94%  '
Set-TextValue 'D35' '11.51'
Set-TextValue 'E35' '  +6.08%  '
Set-TextValue 'D36' '0.02787'
Set-TextValue 'E36' '  +0.00%  '
Set-TextValue 'D37' '0.2514'
Set-TextValue 'E37' '  -0.50%  '
Set-TextValue 'D38' '0.08804'
Set-TextValue 'E38' '  +0.62%  '
Set-TextValue 'D39' '6.076'
Set-TextValue 'E39' '  +0.62%  '
Set-TextValue 'D40' '0.07011'
Set-TextValue 'E40' '  -1.72%  '
Set-TextValue 'E41' '  +4.60%  '
Set-TextValue 'D42' '0.7053'
Set-TextValue 'E42' '  -0.01%  '
Set-TextValue 'D43' '1.328'
Set-TextValue 'E43' '  -1.73%  '
Set-TextValue 'D44' '15.83'
Set-TextValue 'E44' '  -2.79%  '
Set-TextValue 'D45' '0.6504'
Set-TextValue 'E45' '  -0.87%  '
Set-TextValue 'D46' '2.335'
Set-TextValue 'E46' '  +0.44%  '
Set-TextValue 'D47' '0.9965'
Set-TextValue 'E47' '  -0.31%  '
Set-TextValue 'D48' '3.975'
Set-TextValue 'E48' '  -0.11%  '
Set-TextValue 'D49' '0.07973'
Set-TextValue 'E49' '  -0.57%  '
Set-TextValue 'D50' '127.89'
Set-TextValue 'E50' '  +1.40%  '
Set-TextValue 'E51' '  -0.97%  '
